$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 1 (A1:H1) from -0.05 to -2.5
$ws.Range("A1:H1").Value = -2.5

# Update rows 9-14 (A9:H14) from -0.05 to -2.5
$ws.Range("A9:H14").Value = -2.5

# Update the active selection/cell to I12
$ws.Range("I12").Select()
